$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0
$batsman = "Manish Pandey" + $nbsp

# Final header row (ownTeam/oppTeam inserted before batsman; everything after shifts right)
$headers = @("venue", "date", "result", "ownTeam", "oppTeam", "batsman", "totalRuns", "totalBalls", "total4s", "total6s", "sr")

# Final data rows, in final top-to-bottom order:
# venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls, total4s, total6s, sr
$rows = @(
    @(" Dubai (DSC)", " October 27 2020", "Sunrisers won by 88 runs", "Sunrisers Hyderabad", "Delhi Capitals", $batsman, "44", "31", "4", "1", "141.93"),
    @(" Dubai (DSC)", " October 24 2020", "Kings XI won by 12 runs", "Sunrisers Hyderabad", "Kings XI Punjab", $batsman, "15", "29", "0", "0", "51.72"),
    @(" Dubai (DSC)", " October 13 2020", "Super Kings won by 20 runs", "Sunrisers Hyderabad", "Chennai Super Kings", $batsman, "4", "3", "1", "0", "133.33"),
    @(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Sunrisers Hyderabad", "Kolkata Knight Riders", $batsman, "6", "7", "0", "0", "85.71"),
    @(" Dubai (DSC)", " September 21 2020", "RCB won by 10 runs", "Sunrisers Hyderabad", "Royal Challengers Bangalore", $batsman, "34", "33", "3", "1", "103.03"),
    @(" Abu Dhabi", " November 06 2020", "Sunrisers won by 6 wickets (with 2 balls remaining)", "Sunrisers Hyderabad", "Royal Challengers Bangalore", $batsman, "24", "21", "3", "1", "114.28"),
    @(" Dubai (DSC)", " October 22 2020", "Sunrisers won by 8 wickets (with 11 balls remaining)", "Sunrisers Hyderabad", "Rajasthan Royals", $batsman, "83", "47", "4", "8", "176.59"),
    @(" Sharjah", " October 31 2020", "Sunrisers won by 5 wickets (with 35 balls remaining)", "Sunrisers Hyderabad", "Royal Challengers Bangalore", $batsman, "26", "19", "3", "1", "136.84"),
    @(" Abu Dhabi", " September 26 2020", "KKR won by 7 wickets (with 12 balls remaining)", "Sunrisers Hyderabad", "Kolkata Knight Riders", $batsman, "51", "38", "3", "2", "134.21"),
    @(" Abu Dhabi", " November 08 2020", "Capitals won by 17 runs", "Sunrisers Hyderabad", "Delhi Capitals", $batsman, "21", "14", "3", "0", "150.00"),
    @(" Dubai (DSC)", " October 11 2020", "Royals won by 5 wickets (with 1 ball remaining)", "Sunrisers Hyderabad", "Rajasthan Royals", $batsman, "54", "44", "2", "3", "122.72"),
    @(" Dubai (DSC)", " October 08 2020", "Sunrisers won by 69 runs", "Sunrisers Hyderabad", "Kings XI Punjab", $batsman, "1", "2", "0", "0", "50.00"),
    @(" Dubai (DSC)", " October 02 2020", "Sunrisers won by 7 runs", "Sunrisers Hyderabad", "Chennai Super Kings", $batsman, "29", "21", "5", "0", "138.09"),
    @(" Sharjah", " October 04 2020", "Mumbai won by 34 runs", "Sunrisers Hyderabad", "Mumbai Indians", $batsman, "30", "19", "4", "1", "157.89"),
    @(" Abu Dhabi", " September 29 2020", "Sunrisers won by 15 runs", "Sunrisers Hyderabad", "Delhi Capitals", $batsman, "3", "5", "0", "0", "60.00")
)

$totalCols = $headers.Length
$totalRows = $rows.Length + 1

# Format the whole target range as Text first so numeric-looking strings
# ("44", "150.00", ...) are stored as text (matching the source data), not
# auto-converted to numbers.
$fullRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($totalRows, $totalCols))
$fullRange.NumberFormat = "@"

for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}
